$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Insert()

$ws.Range("C2").Value = "https://www.infineon.com/cms/en/product/evaluation-boards/cy8ckit-062s2-43012 "
$ws.Range("A2").Value = "PSoC 6 Wi-Fi/BLE kit"
$ws.Range("B2").Value = 122.43

$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.infineon.com/cms/en/product/evaluation-boards/cy8ckit-062s2-43012 ")
$ws.Range("C2").Style = $ws.Range("C3").Style

$ws.Range("B8").Formula = "=SUM(B2:B6)"

$ws.Range("B9").Select()

Write-Host "done"
